$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7165
$ws.Range("C2").Value = 247
$ws.Range("C3").Value = 149380
$ws.Range("C4").Value = 141101
$ws.Range("C5").Value = 8280
$ws.Range("C6").Value = 900
$ws.Range("C7").Value = 5.54
$ws.Range("C8").Value = 63.47
